$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Power Analysis Tool")

# Update the Minimum Detectable Effect (MDE) input value from 0.04 to 0.05
$ws.Range("C4").Value = 0.05

# Update the formula in C27 (Number of Analysis Units Required) to the new rounding logic
$ws.Range("C27").Formula = "=ROUND(VLOOKUP(`$C`$18,`$B`$21:`$C`$25,2,FALSE)/IF(C6 = 0.5, 1, 1 - ABS((C6 - 0.5)/0.5)),0)"

# Update the active selection on the sheet to G21
$ws.Range("G21").Select()
